$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8 (year 2025) metrics per updated data source
$ws.Range("C8").Value = 1291
$ws.Range("E8").Value = 1085
$ws.Range("G8").Value = 84.04337722695585
$ws.Range("H8").Value = 15.956622773044153
